$d = $word.ActiveDocument

# 1) Title paragraph: merge "Gasoline Price " + "Model" -> "Gasoline Price Model"
$d.Content.Find.Execute(
    "Gasoline Price Model",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "Gasoline Price Model", 2) | Out-Null

# 2) "To consumers..." paragraph: merge the two runs that split right before
#    "a condensed version of the code can be found at "
$oldText2 = "To consumers, perhaps the most frequent price-sensitive purchase is gasoline. Differences of a penny or two can draw consumers from one provider to another. I wondered if I could find data to create a simplistic, first-look model of gasoline prices over the past decade. I used R for this fun exercise; a condensed version of the code can be found at "
$d.Content.Find.Execute(
    $oldText2,
    $false, $false, $false, $false, $false,
    $true, 1, $false, $oldText2, 2) | Out-Null

# 3) Correlation-analysis paragraph: fix typo ("lanned ot include were not
#    hhighly correlated" -> "planned to include were not highly correlated")
#    and merge the many runs into one.
$oldText3 = "I completed some univariate correlation analyses by PADD to verify that the variables I had selected were highly correlated with the price of gasoline and to be sure variables I lanned ot include were not hhighly correlated with each other. The cost of crude oil was highly correlated with the price of gasoline (a Pearson" + [char]0x2019 + "s correlation coefficient [PCC] in the range of 0.95 with very small confidence intervals); federal excise tax on gasoline was not correlated with the price of gasoline as the tax has remained constant at 18.4" + [char]0x00A2 + " per gallon over the time period of interest; state tax and OPEC imports were not correlated with the price of gasoline as the PCC confidence intervals included zero with p values that were non-significant; the summer variable had low but significant PCC values; the year had significant PCC values. The year appears to be an important proxy for shock events, as without the year variable in the model, it performed worse in 2020 and 2021, years affected by the COVID-19 pandemic. Other shocks include unexpected refinery shutdowns due to mechanical issues or natural events. In general, for each PADD, the model underestimated prices."
$newText3 = "I completed some univariate correlation analyses by PADD to verify that the variables I had selected were highly correlated with the price of gasoline and to be sure variables I planned to include were not highly correlated with each other. The cost of crude oil was highly correlated with the price of gasoline (a Pearson" + [char]0x2019 + "s correlation coefficient [PCC] in the range of 0.95 with very small confidence intervals); federal excise tax on gasoline was not correlated with the price of gasoline as the tax has remained constant at 18.4" + [char]0x00A2 + " per gallon over the time period of interest; state tax and OPEC imports were not correlated with the price of gasoline as the PCC confidence intervals included zero with p values that were non-significant; the summer variable had low but significant PCC values; the year had significant PCC values. The year appears to be an important proxy for shock events, as without the year variable in the model, it performed worse in 2020 and 2021, years affected by the COVID-19 pandemic. Other shocks include unexpected refinery shutdowns due to mechanical issues or natural events. In general, for each PADD, the model underestimated prices."
$d.Content.Find.Execute(
    $oldText3,
    $false, $false, $false, $false, $false,
    $true, 1, $false, $newText3, 2) | Out-Null

# 4) "In general, the fit of the models..." paragraph: merge the three runs
#    into one (no text change).
$oldText4 = "In general, the fit of the models was acceptable, with the PADD1, PADD2 and PADD3 models having the best fit, as expected from the ANOVA results above. The adjusted r-squared values for the models were: 0.9154, 0.9139, 0.9360, 0.8202 and 0.8201 for PADDs one through five respectively and the sum of the residuals was zero in all cases."
$d.Content.Find.Execute(
    $oldText4,
    $false, $false, $false, $false, $false,
    $true, 1, $false, $oldText4, 2) | Out-Null

# 5) "The variance across models..." paragraph: reword + merge runs into one.
$oldText5 = "The variance across models of the coefficients for the cost of crude oil is relatively small. The values show that for each increase of `$1 in the cost of crude oil, the price of gasoline can be expected to increase, on average across the PADDs, around 3.13" + [char]0x00A2 + ". The coefficients for the summer variable show the most variance across the five models, with switching to summer blend costing an additional 10.12" + [char]0x00A2 + " in PADD1 (East Coast) to 28.68" + [char]0x00A2 + " in PADD4 (Rocky Mountain Region)."
$newText5 = "The variance of the coefficients across models for the cost of crude oil is relatively small. The values show that for each increase of `$1 in the cost of crude oil, the price of gasoline can be expected to increase, on average across the PADDs, around 3.13" + [char]0x00A2 + ". The coefficients for the summer variable show the most variance across the five models, with switching to summer blend costing an additional 10.12" + [char]0x00A2 + " in PADD1 (East Coast) to 28.68" + [char]0x00A2 + " in PADD4 (Rocky Mountain Region)."
$d.Content.Find.Execute(
    $oldText5,
    $false, $false, $false, $false, $false,
    $true, 1, $false, $newText5, 2) | Out-Null

Write-Host "Done."
